# "N de BdL" sheet: columns D (rows 2-11) and the stray "0 pal" cells in
# column E (rows 6-9) were stored as useless text-with-a-comma strings
# ("0" / " 0,") instead of the plain number 0. Reformat those cells as
# numbers and overwrite them with a literal 0 so they become real numeric
# cells and the bogus strings disappear from the shared-string table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D11").NumberFormat = "general"
$ws.Range("D2:D11").Value = 0

$ws.Range("E6:E9").NumberFormat = "general"
$ws.Range("E6:E9").Value = 0

[void]$ws.Range("E13").Select()
